# Auto-generated Excel COM-interop script
# Applies weekly crime-data refresh (new week's figures) to cs-en-us-121pct.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -------------------------------------------------
$ws.Range("A8").Characters(21, 2).Text = "44"
$ws.Range("C9").Characters(27, 10).Text = "10/31/2022"
$ws.Range("C9").Characters(48, 10).Text = "11/6/2022"

# --- Data table updates ---------------------------------------------------
# Row 15
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E15").Value = -100
$ws.Range("F15").NumberFormat = '#,##0'
$ws.Range("F15").Value = 1
$ws.Range("G15").NumberFormat = '#,##0'
$ws.Range("G15").Value = 2
$ws.Range("H15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H15").Value = -50
$ws.Range("J15").NumberFormat = '#,##0'
$ws.Range("J15").Value = 16
$ws.Range("K15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K15").Value = -37.5

# Row 16
$ws.Range("C16").NumberFormat = '#,##0'
$ws.Range("C16").Value = 1
$ws.Range("E16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E16").Value = -50
$ws.Range("G16").NumberFormat = '#,##0'
$ws.Range("G16").Value = 8
$ws.Range("H16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H16").Value = 25
$ws.Range("I16").NumberFormat = '#,##0'
$ws.Range("I16").Value = 67
$ws.Range("J16").NumberFormat = '#,##0'
$ws.Range("J16").Value = 57
$ws.Range("K16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K16").Value = 17.543859649122
$ws.Range("L16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L16").Value = 6.349206349206

# Row 17
$ws.Range("C17").NumberFormat = '#,##0'
$ws.Range("C17").Value = 1
$ws.Range("D17").NumberFormat = '#,##0'
$ws.Range("D17").Value = 5
$ws.Range("E17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E17").Value = -80
$ws.Range("F17").NumberFormat = '#,##0'
$ws.Range("F17").Value = 12
$ws.Range("G17").NumberFormat = '#,##0'
$ws.Range("G17").Value = 13
$ws.Range("H17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H17").Value = -7.692307692307
$ws.Range("I17").NumberFormat = '#,##0'
$ws.Range("I17").Value = 183
$ws.Range("J17").NumberFormat = '#,##0'
$ws.Range("J17").Value = 138
$ws.Range("K17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K17").Value = 32.608695652173
$ws.Range("L17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L17").Value = 13.664596273291

# Row 18
$ws.Range("D18").NumberFormat = '#,##0'
$ws.Range("D18").Value = 1
$ws.Range("F18").NumberFormat = '#,##0'
$ws.Range("F18").Value = 4
$ws.Range("G18").NumberFormat = '#,##0'
$ws.Range("G18").Value = 7
$ws.Range("H18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H18").Value = -42.857142857142
$ws.Range("J18").NumberFormat = '#,##0'
$ws.Range("J18").Value = 72
$ws.Range("K18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K18").Value = 1.388888888888

# Row 19
$ws.Range("C19").NumberFormat = '#,##0'
$ws.Range("C19").Value = 11
$ws.Range("D19").NumberFormat = '#,##0'
$ws.Range("D19").Value = 7
$ws.Range("E19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E19").Value = 57.142857142857
$ws.Range("F19").NumberFormat = '#,##0'
$ws.Range("F19").Value = 33
$ws.Range("G19").NumberFormat = '#,##0'
$ws.Range("G19").Value = 25
$ws.Range("H19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H19").Value = 32
$ws.Range("I19").NumberFormat = '#,##0'
$ws.Range("I19").Value = 346
$ws.Range("J19").NumberFormat = '#,##0'
$ws.Range("J19").Value = 268
$ws.Range("K19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K19").Value = 29.10447761194
$ws.Range("L19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L19").Value = 45.378151260504

# Row 20
$ws.Range("C20").NumberFormat = '#,##0'
$ws.Range("C20").Value = 3
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Value = "'0"
$ws.Range("E20").NumberFormat = "General"
$ws.Range("E20").Value = "'***.*"
$ws.Range("F20").NumberFormat = '#,##0'
$ws.Range("F20").Value = 14
$ws.Range("H20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H20").Value = 180
$ws.Range("I20").NumberFormat = '#,##0'
$ws.Range("I20").Value = 91
$ws.Range("K20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K20").Value = 111.627906976744
$ws.Range("L20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L20").Value = 102.222222222222

# Row 21
$ws.Range("C21").NumberFormat = '#,##0'
$ws.Range("C21").Value = 16
$ws.Range("D21").NumberFormat = '#,##0'
$ws.Range("D21").Value = 16
$ws.Range("E21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("E21").Value = 0
$ws.Range("G21").NumberFormat = '#,##0'
$ws.Range("G21").Value = 60
$ws.Range("H21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("H21").Value = 23.333333333333
$ws.Range("I21").NumberFormat = '#,##0'
$ws.Range("I21").Value = 770
$ws.Range("J21").NumberFormat = '#,##0'
$ws.Range("J21").Value = 597
$ws.Range("K21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("K21").Value = 28.978224455611
$ws.Range("L21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("L21").Value = 27.906976744186

# Row 23
$ws.Range("G23").NumberFormat = '#,##0'
$ws.Range("G23").Value = 2
$ws.Range("H23").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H23").Value = -50
$ws.Range("J23").NumberFormat = '#,##0'
$ws.Range("J23").Value = 13
$ws.Range("K23").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K23").Value = 30.76923076923

# Row 24
$ws.Range("C24").NumberFormat = '#,##0'
$ws.Range("C24").Value = 27
$ws.Range("D24").NumberFormat = '#,##0'
$ws.Range("D24").Value = 19
$ws.Range("E24").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E24").Value = 42.105263157894
$ws.Range("F24").NumberFormat = '#,##0'
$ws.Range("F24").Value = 97
$ws.Range("G24").NumberFormat = '#,##0'
$ws.Range("G24").Value = 66
$ws.Range("H24").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H24").Value = 46.969696969697
$ws.Range("I24").NumberFormat = '#,##0'
$ws.Range("I24").Value = 1192
$ws.Range("J24").NumberFormat = '#,##0'
$ws.Range("J24").Value = 781
$ws.Range("K24").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K24").Value = 52.624839948783
$ws.Range("L24").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L24").Value = 66.480446927374

# Row 25
$ws.Range("C25").NumberFormat = '#,##0'
$ws.Range("C25").Value = 8
$ws.Range("D25").NumberFormat = '#,##0'
$ws.Range("D25").Value = 9
$ws.Range("E25").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E25").Value = -11.111111111111
$ws.Range("F25").NumberFormat = '#,##0'
$ws.Range("F25").Value = 42
$ws.Range("G25").NumberFormat = '#,##0'
$ws.Range("G25").Value = 33
$ws.Range("H25").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H25").Value = 27.272727272727
$ws.Range("I25").NumberFormat = '#,##0'
$ws.Range("I25").Value = 424
$ws.Range("J25").NumberFormat = '#,##0'
$ws.Range("J25").Value = 346
$ws.Range("K25").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K25").Value = 22.543352601156
$ws.Range("L25").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L25").Value = 39.473684210526

# Row 26
$ws.Range("F26").NumberFormat = '#,##0'
$ws.Range("F26").Value = 1
$ws.Range("G26").NumberFormat = '#,##0'
$ws.Range("G26").Value = 5
$ws.Range("H26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H26").Value = -80
$ws.Range("I26").NumberFormat = '#,##0'
$ws.Range("I26").Value = 16
$ws.Range("J26").NumberFormat = '#,##0'
$ws.Range("J26").Value = 28
$ws.Range("K26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K26").Value = -42.857142857142
$ws.Range("L26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L26").Value = -23.809523809523

# Row 27
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("D27").Value = 3
$ws.Range("J27").NumberFormat = '#,##0'
$ws.Range("J27").Value = 41
$ws.Range("K27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K27").Value = -7.317073170731
$ws.Range("L27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L27").Value = 72.727272727272

# Row 28
$ws.Range("G28").NumberFormat = "General"
$ws.Range("G28").Value = "'0"
$ws.Range("H28").NumberFormat = "General"
$ws.Range("H28").Value = "'***.*"

# Row 29
$ws.Range("G29").NumberFormat = "General"
$ws.Range("G29").Value = "'0"
$ws.Range("H29").NumberFormat = "General"
$ws.Range("H29").Value = "'***.*"

